# Applies the "Diff and initial highlighting now work" edit:
#  - Sheet1 gains a new "Origin" column between Fruit and Color, and the
#    leading note text changes; the trailing junk cell is cleared.
#  - Sheet2's table is rebuilt: the Date column is dropped, a new "Shape"
#    column is added, the sample rows are replaced, and the stray trailing
#    row is removed.
#  - The active sheet/tab flips from Sheet1 to Sheet2, and each sheet's
#    remembered selection moves to A5.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# ---------------------------------------------------------------------
# Sheet1: insert "Origin" column (D), shifting the old Color column to E.
# ---------------------------------------------------------------------
$ws1.Columns.Item(4).Insert()

$ws1.Range("D4").Value = "Origin"
$ws1.Range("D5").Value = "Maine"
$ws1.Range("D6").Value = "Chile"
$ws1.Range("D7").Value = "Hawaii"
$ws1.Range("D8").Value = "Thailand"

# The old junk marker ("END") shifted from E8 into F8 -- drop it, leaving
# an empty, styled cell behind (matches F9/F10's blank marker cells).
$ws1.Range("F8").ClearContents()

# Leading note text above the table.
$ws1.Range("A2").Value = "This is a leading line, before the table"

# ---------------------------------------------------------------------
# Sheet2: rebuild the table -- drop Date, add Shape, replace sample data.
# ---------------------------------------------------------------------
$ws2.Range("A2").Value = "This is a leading line, before the table"

$ws2.Range("A4").Value = "ID"
$ws2.Range("B4").Value = "Fruit"
$ws2.Range("C4").Value = "Shape"
$ws2.Range("D4").Value = "Color"

$ws2.Range("A5:D7").ClearContents()
$ws2.Range("A5").Value = 5
$ws2.Range("B5").Value = "grapes"
$ws2.Range("C5").Value = "pebbles"
$ws2.Range("D5").Value = "purple"

$ws2.Range("A6").Value = 2
$ws2.Range("B6").Value = "mango"
$ws2.Range("C6").Value = "oval"
$ws2.Range("D6").Value = "yellow"

$ws2.Range("A7").Value = 4
$ws2.Range("B7").Value = "dried mango"
$ws2.Range("C7").Value = "flat"
$ws2.Range("D7").Value = "orange"

$ws2.Range("A5:A7").NumberFormat = "General"

# Row 8 no longer carries data -- delete the (now unused) A:D cells
# outright so only the E8 blank marker cell remains, matching E9's
# pattern, then remove the now-unused trailing row 10.
$ws2.Range("A8:D8").Delete()
$ws2.Range("E8").ClearContents()
$ws2.Rows.Item(10).Delete()

# ---------------------------------------------------------------------
# Selection / active-tab bookkeeping.
# ---------------------------------------------------------------------
$ws1.Range("A5").Select()
$ws2.Range("A5").Select()

$ws2.Activate()
